$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A54").Value = "'10/25/2025"
$ws.Range("A54").Style = "Normal"
$ws.Range("B54").Value = 0.190092683663864
$ws.Range("C54").Value = 0.809907316336136
